$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.241.33'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '1.905.83'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '307.39'
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D7").Value = '0.5269'
$ws.Range("E7").Value = '  +1.16%  '

$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  +1.34%  '

$ws.Range("D9").Value = '0.07284'
$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("D10").Value = '21.79'
$ws.Range("E10").Value = '  +2.83%  '

$ws.Range("D11").Value = '0.9031'
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").Value = '0.08181'
$ws.Range("E12").Value = '  -3.33%  '

$ws.Range("D13").Value = '96.31'
$ws.Range("E13").Value = '  -0.66%  '

$ws.Range("D14").Value = '5.364'
$ws.Range("E14").Value = '  +1.23%  '

$ws.Range("D15").Value = '1.455.25'
$ws.Range("E15").Value = '  -23.72%  '

$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").Value = '0.000008655'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("E19").Value = '  -0.10%  '

$ws.Range("D20").Value = '27.281.03'
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").Value = '6.516'
$ws.Range("E23").Value = '  +1.18%  '

$ws.Range("D24").Value = '150.07'
$ws.Range("E24").Value = '  +2.07%  '

$ws.Range("D25").Value = '2.311'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -0.69%  '

$ws.Range("D28").Value = "'116.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.30%  '

$ws.Range("D29").Value = '4.847'
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("D30").Value = '4.854'
$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("D31").Value = '0.09246'
$ws.Range("E31").Value = '  -0.45%  '

$ws.Range("D32").Value = '0.8299'
$ws.Range("E32").Value = '  +4.29%  '

$ws.Range("D33").Value = '0.05059'
$ws.Range("E33").Value = '  -0.18%  '

$ws.Range("E34").Value = '  -0.75%  '

$ws.Range("D35").Value = '2.993'
$ws.Range("E35").Value = '  +1.46%  '

# Rows 36-37: reorder (MXToken now above RenderToken)
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = '3.352'
$ws.Range("E36").Value = '  -2.54%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").Value = '  +5.56%  '

$ws.Range("D38").Value = '0.5809'
$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("D39").Value = '0.02004'
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").Value = '1.082'
$ws.Range("E40").Value = '  +0.55%  '

$ws.Range("D41").Value = "'9.170"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.26%  '

$ws.Range("D42").Value = '6.613'
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").Value = '117.33'
$ws.Range("E43").Value = '  +0.89%  '

$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("E45").Value = '  +1.34%  '

# Rows 46-47: reorder (EnergySwap now above PaxDollar)
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.21'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").Value = '1.643'
$ws.Range("E48").Value = '  +0.60%  '

$ws.Range("D49").Value = '38.99'
$ws.Range("E49").Value = '  +3.46%  '

$ws.Range("D50").Value = '0.06144'
$ws.Range("E50").Value = '  +3.03%  '

$ws.Range("D51").Value = '64.49'
$ws.Range("E51").Value = '  +0.75%  '
